# Rotates the species-record data among rows 3, 4 and 6:
#   row3 <- old row4 data, row4 <- old row6 data, row6 <- old row3 data
# (columns A, B, E, F, G, H, Q, R, AC - all other columns already match
# across these three rows, so they are left untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: becomes the old Row 4 content ---
$ws.Range("A3").Value = 111638283
$ws.Range("B3").Value = 77515
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 588212.5872669512
$ws.Range("R3").Value = 7033297.86989607
$ws.Range("AC3").ClearContents()

# --- Row 4: becomes the old Row 6 content ---
$ws.Range("A4").Value = 111638282
$ws.Range("B4").Value = 89405
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 588265.5989708689
$ws.Range("R4").Value = 7033311.880202802

# --- Row 6: becomes the old Row 3 content ---
$ws.Range("A6").Value = 111638278
$ws.Range("B6").Value = 56414
$ws.Range("E6").Value = 100049
$ws.Range("F6").Value = "Spillkråka"
$ws.Range("G6").Value = "Dryocopus martius"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("Q6").Value = 588319.7233830886
$ws.Range("R6").Value = 7033285.591169797
$ws.Range("AC6").Value = "Bohål i gammal grov tall."
